$wb = $excel.ActiveWorkbook

# Clone the last existing sheet ("ODI Bowling") to the end of the workbook so
# the new sheet inherits the same sheetPr / pageMargins / header style as the
# rest of the workbook, then rename it and replace its contents.
$sourceSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sourceSheet.Copy($null, $sourceSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "ODI Batting Extra"

# Drop the inherited values/format so only the columns we need remain.
$ws.Cells.ClearContents()
$ws.Columns.Item(7).Clear()

# Header row (keeps the bold/centered/bordered style copied from the source sheet)
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Data row: MATCH_CODE / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH are
# stored as text, BATTING_POSITION is a genuine number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "4679"
$ws.Range("B2").Value = 11
$ws.Range("C2:D2").NumberFormat = "@"
$ws.Range("C2").Value = "0"
$ws.Range("D2").Value = "0"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.08%"
$ws.Range("F2").Value = "NO"

# Leave the original first sheet as the active/selected one, as in the source file.
$wb.Worksheets.Item(1).Activate()
